$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'42.065.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.34%  "

# Row 3
$ws.Range("D3").Value = "'2.218.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.17%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'241.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.25%  "

# Row 6
$ws.Range("D6").Value = "'0.629"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.11%  "

# Row 7
$ws.Range("D7").Value = "'73.29"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.42%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("E9").Value = "  -1.66%  "

# Row 10
$ws.Range("D10").Value = "'42.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.87%  "

# Row 11
$ws.Range("D11").Value = "'0.0955"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.94%  "

# Row 12
$ws.Range("D12").Value = "'7.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.23%  "

# Row 13
$ws.Range("E13").Value = "  +0.01%  "

# Row 14
$ws.Range("D14").Value = "'2.549.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.08%  "

# Row 15
$ws.Range("D15").Value = "'14.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.19%  "

# Row 16
$ws.Range("D16").Value = "'0.838"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.83%  "

# Row 17
$ws.Range("D17").Value = "'2.224.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.85%  "

# Row 18
$ws.Range("D18").Value = "'41.897.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.48%  "

# Row 19
$ws.Range("D19").Value = "'0.0000108"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.96%  "

# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.55%  "

# Row 21
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Value = "'72.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.73%  "

# Row 22
$ws.Range("D22").Value = "'10.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +17.00%  "

# Row 23
$ws.Range("D23").Value = "'229.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.55%  "

# Row 24
$ws.Range("D24").Value = "'2.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.99%  "

# Row 25
$ws.Range("D25").Value = "'11.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.86%  "

# Row 26
$ws.Range("E26").Value = "  +0.08%  "

# Row 27
$ws.Range("D27").Value = "'3.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "

# Row 28
$ws.Range("E28").Value = "  -1.48%  "

# Row 29
$ws.Range("D29").Value = "'2.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.63%  "

# Row 30
$ws.Range("D30").Value = "'167.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.88%  "

# Row 31
$ws.Range("D31").Value = "'20.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.82%  "

# Row 32
$ws.Range("D32").Value = "'5.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.71%  "

# Row 33
$ws.Range("D33").Value = "'0.0799"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.70%  "

# Row 34
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "'0.125"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.44%  "

# Row 35
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'29.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.61%  "

# Row 36
$ws.Range("E36").Value = "  -9.26%  "

# Row 37
$ws.Range("D37").Value = "'4.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.88%  "

# Row 38
$ws.Range("E38").Value = "  -4.29%  "

# Row 39
$ws.Range("D39").Value = "'13.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.01%  "

# Row 40
$ws.Range("D40").Value = "'65.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.81%  "

# Row 41
$ws.Range("E41").Value = "  -2.28%  "

# Row 42
$ws.Range("E42").Value = "  -2.54%  "

# Row 43
$ws.Range("E43").Value = "  -2.80%  "

# Row 44
$ws.Range("D44").Value = "'8.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.26%  "

# Row 45
$ws.Range("D45").Value = "'105.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.53%  "

# Row 46
$ws.Range("E46").Value = "  -2.38%  "

# Row 47
$ws.Range("E47").Value = "  +4.34%  "

# Row 48
$ws.Range("E48").Value = "  +0.30%  "

# Row 49
$ws.Range("E49").Value = "  -0.68%  "

# Row 50
$ws.Range("E50").Value = "  -0.10%  "

# Row 51
$ws.Range("D51").Value = "'2.425.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.22%  "
